$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row at position 8 (shifts the existing data rows,
#    the totals row, and the footer row down by one).
# ------------------------------------------------------------------
$ws.Rows("8:8").Insert()

# ------------------------------------------------------------------
# 2. The newly inserted row 8 comes back with a blank/default style;
#    restore the same per-column formatting used by every other data
#    row by copying the formats from row 9 (which now holds what used
#    to be row 8's formatting) and fix the row height.
# ------------------------------------------------------------------
$ws.Range("A9:Q9").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$ws.Rows("8:8").RowHeight = 24.75

# ------------------------------------------------------------------
# 3. Fill in the new item's data: CONTAFEVER N 200MG/5ML SUSP. 120ML
# ------------------------------------------------------------------
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "CONTAFEVER N 200MG/5ML SUSP. 120ML"
$ws.Range("H8").Value = "9:0"
$ws.Range("N8").Value = "33.00"
$ws.Range("Q8").Value = "1:0"

# L8 and P8 sit on cells whose number format is numeric, but the
# workbook stores these figures as literal text - force text storage
# with a leading quote, then re-apply the clean (non quote-prefixed)
# format from a sibling cell so the style index matches the rest of
# the table exactly.
$ws.Range("L8").Value = "'1"
$ws.Range("L7").Copy()
$ws.Range("L8").PasteSpecial(-4122)

$ws.Range("P8").Value = "'33.0000"
$ws.Range("P7").Copy()
$ws.Range("P8").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4. The totals row (old row 15) is now row 16; its height grows to
#    accommodate the extra item and its value must include the new
#    item's selling price.
# ------------------------------------------------------------------
$ws.Rows("16:16").RowHeight = 25.5
$ws.Range("P16").Value = 489.8

# ------------------------------------------------------------------
# 5. The footer row (old row 16) is now row 17; refresh the printed
#    timestamp to reflect the new export time.
# ------------------------------------------------------------------
$ws.Range("A17").Value = "Sunday, 5 October, 2025 10:15 AM"
